$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2928123333333333
$ws.Range("H2").Value = 0.8784369999999999
$ws.Range("I2").Value = 0.06406943071632207
$ws.Range("J2").Value = 0.06406943071632207
$ws.Range("M2").Value = 6.391557333333332
$ws.Range("N2").Value = 19.174672
$ws.Range("O2").Value = 0.1156448793857254
$ws.Range("P2").Value = 0.1156448793857254
$ws.Range("Q2").Value = 1.87152681640711
$ws.Range("R2").Value = 16.843741347664
$ws.Range("S2").Value = 0.007409301587501154
$ws.Range("T2").Value = 0.007409301587501154

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2928123333333333
$ws.Range("H3").Value = 0.8784369999999999
$ws.Range("I3").Value = 0.06406943071632207
$ws.Range("J3").Value = 0.06406943071632207
$ws.Range("O3").Value = 0.2610362896883882
$ws.Range("P3").Value = 0.2610362896883882
$ws.Range("Q3").Value = 4.224453506304888
$ws.Range("R3").Value = 38.020081556744
$ws.Range("S3").Value = 0.01672444647663596
$ws.Range("T3").Value = 0.01672444647663596

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2928123333333333
$ws.Range("H4").Value = 0.8784369999999999
$ws.Range("I4").Value = 0.06406943071632207
$ws.Range("J4").Value = 0.06406943071632207
$ws.Range("M4").Value = 32.348972
$ws.Range("N4").Value = 97.046916
$ws.Range("O4").Value = 0.5853022620452971
$ws.Range("P4").Value = 0.5853022620452972
$ws.Range("Q4").Value = 9.472177972254665
$ws.Range("R4").Value = 85.24960175029199
$ws.Range("S4").Value = 0.03749998272621775
$ws.Range("T4").Value = 0.03749998272621775

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2928123333333333
$ws.Range("H5").Value = 0.8784369999999999
$ws.Range("I5").Value = 0.06406943071632207
$ws.Range("J5").Value = 0.06406943071632207
$ws.Range("M5").Value = 2.101131333333333
$ws.Range("N5").Value = 6.303394
$ws.Range("O5").Value = 0.03801656888058921
$ws.Range("P5").Value = 0.03801656888058921
$ws.Range("Q5").Value = 0.615237168353111
$ws.Range("R5").Value = 5.537134515178
$ws.Range("S5").Value = 0.002435699925967196
$ws.Range("T5").Value = 0.002435699925967196

$ws.Range("I6").Value = 0.8630927339690215
$ws.Range("J6").Value = 0.8630927339690215
$ws.Range("M6").Value = 6.391557333333332
$ws.Range("N6").Value = 19.174672
$ws.Range("O6").Value = 0.1156448793857254
$ws.Range("P6").Value = 0.1156448793857254
$ws.Range("Q6").Value = 25.21173012791644
$ws.Range("R6").Value = 226.905571151248
$ws.Range("S6").Value = 0.09981225511854346
$ws.Range("T6").Value = 0.09981225511854346

$ws.Range("I7").Value = 0.8630927339690215
$ws.Range("J7").Value = 0.8630927339690215
$ws.Range("O7").Value = 0.2610362896883882
$ws.Range("P7").Value = 0.2610362896883882
$ws.Range("S7").Value = 0.2252985249322805
$ws.Range("T7").Value = 0.2252985249322805

$ws.Range("I8").Value = 0.8630927339690215
$ws.Range("J8").Value = 0.8630927339690215
$ws.Range("M8").Value = 32.348972
$ws.Range("N8").Value = 97.046916
$ws.Range("O8").Value = 0.5853022620452971
$ws.Range("P8").Value = 0.5853022620452972
$ws.Range("Q8").Value = 127.6016953999826
$ws.Range("R8").Value = 1148.415258599844
$ws.Range("S8").Value = 0.5051701295469282
$ws.Range("T8").Value = 0.5051701295469282

$ws.Range("I9").Value = 0.8630927339690215
$ws.Range("J9").Value = 0.8630927339690215
$ws.Range("M9").Value = 2.101131333333333
$ws.Range("N9").Value = 6.303394
$ws.Range("O9").Value = 0.03801656888058921
$ws.Range("P9").Value = 0.03801656888058921
$ws.Range("Q9").Value = 8.287988885438445
$ws.Range("R9").Value = 74.591899968946
$ws.Range("S9").Value = 0.03281182437126937
$ws.Range("T9").Value = 0.03281182437126937

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.332886
$ws.Range("H10").Value = 0.998658
$ws.Range("I10").Value = 0.07283783531465635
$ws.Range("J10").Value = 0.07283783531465635
$ws.Range("M10").Value = 6.391557333333332
$ws.Range("N10").Value = 19.174672
$ws.Range("O10").Value = 0.1156448793857254
$ws.Range("P10").Value = 0.1156448793857254
$ws.Range("Q10").Value = 2.127659954464
$ws.Range("R10").Value = 19.148939590176
$ws.Range("S10").Value = 0.008423322679680762
$ws.Range("T10").Value = 0.008423322679680762

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.332886
$ws.Range("H11").Value = 0.998658
$ws.Range("I11").Value = 0.07283783531465635
$ws.Range("J11").Value = 0.07283783531465635
$ws.Range("O11").Value = 0.2610362896883882
$ws.Range("P11").Value = 0.2610362896883882
$ws.Range("Q11").Value = 4.802603134544
$ws.Range("R11").Value = 43.223428210896
$ws.Range("S11").Value = 0.01901331827947175
$ws.Range("T11").Value = 0.01901331827947175

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.332886
$ws.Range("H12").Value = 0.998658
$ws.Range("I12").Value = 0.07283783531465635
$ws.Range("J12").Value = 0.07283783531465635
$ws.Range("M12").Value = 32.348972
$ws.Range("N12").Value = 97.046916
$ws.Range("O12").Value = 0.5853022620452971
$ws.Range("P12").Value = 0.5853022620452972
$ws.Range("Q12").Value = 10.768519893192
$ws.Range("R12").Value = 96.916679038728
$ws.Range("S12").Value = 0.04263214977215118
$ws.Range("T12").Value = 0.04263214977215119

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.332886
$ws.Range("H13").Value = 0.998658
$ws.Range("I13").Value = 0.07283783531465635
$ws.Range("J13").Value = 0.07283783531465635
$ws.Range("M13").Value = 2.101131333333333
$ws.Range("N13").Value = 6.303394
$ws.Range("O13").Value = 0.03801656888058921
$ws.Range("P13").Value = 0.03801656888058921
$ws.Range("Q13").Value = 0.6994372050280001
$ws.Range("R13").Value = 6.294934845252
$ws.Range("S13").Value = 0.002769044583352646
$ws.Range("T13").Value = 0.002769044583352646
